# Generate Report for Handback
#
# The localization-status workbook tracks two files:
#   5f2c9c7b-8ee7-48fd-af0c-5866b81aa82d.md
#   0e314636-8642-429d-95e2-56fccc4a9f14.md
#
# 0e314636-... has now been handed back (it previously was only
# "Ready for handoff"). This swaps its row into the "row 2" slot of every
# sheet (matching the other, already-handed-back file's row layout) and
# updates its status / handback timestamps accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A2").Value = "0e314636-8642-429d-95e2-56fccc4a9f14.md"
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("D2").Value = "2016-03-22 08:46:41"

$overview.Range("A3").Value = "5f2c9c7b-8ee7-48fd-af0c-5866b81aa82d.md"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"
$overview.Range("D3").Value = "2016-03-22 08:45:14"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A2").Value = "0e314636-8642-429d-95e2-56fccc4a9f14.md"
$zhcn.Range("B2").Value = ".md"
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("D2").Value = "0e314636-8642-429d-95e2-56fccc4a9f14.2981ce20929d003ce22b02035c8278eea0ddbf86.zh-cn.xlf"
$zhcn.Range("E2").Value = "2016-03-22 08:46:37"
$zhcn.Range("F2").Value = "0e314636-8642-429d-95e2-56fccc4a9f14.md"
$zhcn.Range("G2").Value = "0e314636-8642-429d-95e2-56fccc4a9f14.2981ce20929d003ce22b02035c8278eea0ddbf86.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-03-22 08:47:00"
$zhcn.Range("J2").Value = "Include"

$zhcn.Range("A3").Value = "5f2c9c7b-8ee7-48fd-af0c-5866b81aa82d.md"
$zhcn.Range("B3").Value = ".md"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("D3").Value = "5f2c9c7b-8ee7-48fd-af0c-5866b81aa82d.09e02514cb63c024c396759722709106f98eec03.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-22 08:45:10"
$zhcn.Range("F3").Value = "5f2c9c7b-8ee7-48fd-af0c-5866b81aa82d.md"
$zhcn.Range("G3").Value = "5f2c9c7b-8ee7-48fd-af0c-5866b81aa82d.09e02514cb63c024c396759722709106f98eec03.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-03-22 08:45:49"
$zhcn.Range("J3").Value = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A2").Value = "0e314636-8642-429d-95e2-56fccc4a9f14.md"
$dede.Range("B2").Value = ".md"
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("D2").Value = "0e314636-8642-429d-95e2-56fccc4a9f14.2981ce20929d003ce22b02035c8278eea0ddbf86.de-de.xlf"
$dede.Range("E2").Value = "2016-03-22 08:46:41"
$dede.Range("F2").Value = "0e314636-8642-429d-95e2-56fccc4a9f14.md"
$dede.Range("G2").Value = "0e314636-8642-429d-95e2-56fccc4a9f14.2981ce20929d003ce22b02035c8278eea0ddbf86.de-de.xlf"
$dede.Range("H2").Value = "2016-03-22 08:47:08"
$dede.Range("J2").Value = "Include"

$dede.Range("A3").Value = "5f2c9c7b-8ee7-48fd-af0c-5866b81aa82d.md"
$dede.Range("B3").Value = ".md"
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("D3").Value = "5f2c9c7b-8ee7-48fd-af0c-5866b81aa82d.09e02514cb63c024c396759722709106f98eec03.de-de.xlf"
$dede.Range("E3").Value = "2016-03-22 08:45:14"
$dede.Range("F3").Value = "5f2c9c7b-8ee7-48fd-af0c-5866b81aa82d.md"
$dede.Range("G3").Value = "5f2c9c7b-8ee7-48fd-af0c-5866b81aa82d.09e02514cb63c024c396759722709106f98eec03.de-de.xlf"
$dede.Range("H3").Value = "2016-03-22 08:45:57"
$dede.Range("J3").Value = "Include"
